$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 '61.363.48'
Set-TextCell 2 5 '  +2.11%  '

# Row 3
Set-TextCell 3 4 '2.656.90'
Set-TextCell 3 5 '  +2.56%  '

# Row 4
Set-TextCell 4 4 '1.00'
Set-TextCell 4 5 '  +0.13%  '

# Row 5
Set-TextCell 5 4 '580.96'
Set-TextCell 5 5 '  +0.41%  '

# Row 6
Set-TextCell 6 4 '144.43'
Set-TextCell 6 5 '  +1.51%  '

# Row 7
Set-TextCell 7 4 '0.998'
Set-TextCell 7 5 '  -0.06%  '

# Row 8
Set-TextCell 8 4 '0.599'
Set-TextCell 8 5 '  +0.61%  '

# Row 9
Set-TextCell 9 5 '  +1.43%  '

# Row 10
Set-TextCell 10 5 '  +4.80%  '

# Row 11
Set-TextCell 11 5 '  +3.42%  '

# Row 12
Set-TextCell 12 5 '  +1.15%  '

# Row 13
Set-TextCell 13 4 '3.143.17'
Set-TextCell 13 5 '  +2.98%  '

# Row 14
Set-TextCell 14 4 '26.17'
Set-TextCell 14 5 '  +6.47%  '

# Row 15
Set-TextCell 15 4 '61.331.01'
Set-TextCell 15 5 '  +2.06%  '

# Row 16
Set-TextCell 16 4 '0.0000147'
Set-TextCell 16 5 '  +4.18%  '

# Row 17
Set-TextCell 17 4 '2.665.29'
Set-TextCell 17 5 '  +2.48%  '

# Row 18
Set-TextCell 18 4 '11.70'
Set-TextCell 18 5 '  +2.45%  '

# Row 19
Set-TextCell 19 5 '  +3.51%  '

# Row 20
Set-TextCell 20 4 '354.89'
Set-TextCell 20 5 '  +2.62%  '

# Row 21
Set-TextCell 21 4 '6.90'
Set-TextCell 21 5 '  +0.39%  '

# Row 22
Set-TextCell 22 4 '1.00'
Set-TextCell 22 5 '  +0.18%  '

# Row 23
Set-TextCell 23 4 '0.527'
Set-TextCell 23 5 '  +0.53%  '

# Row 24
Set-TextCell 24 4 '64.51'
Set-TextCell 24 5 '  +2.70%  '

# Row 25
Set-TextCell 25 5 '  +3.51%  '

# Row 26
Set-TextCell 26 5 '  +6.53%  '

# Row 27
Set-TextCell 27 4 '0.998'
Set-TextCell 27 5 '  -0.13%  '

# Row 28
Set-TextCell 28 5 '  +8.00%  '

# Row 29
Set-TextCell 29 5 '  +4.65%  '

# Row 30
Set-TextCell 30 4 '6.94'
Set-TextCell 30 5 '  +9.10%  '

# Row 31
Set-TextCell 31 4 '168.82'
Set-TextCell 31 5 '  +2.76%  '

# Row 32
Set-TextCell 32 5 '  -0.03%  '

# Row 33
Set-TextCell 33 4 '20.18'
Set-TextCell 33 5 '  +4.11%  '

# Row 34
Set-TextCell 34 4 '1.12'
Set-TextCell 34 5 '  +15.39%  '

# Row 35
Set-TextCell 35 5 '  +10.26%  '

# Row 36
Set-TextCell 36 5 '  +11.30%  '

# Row 37
Set-TextCell 37 2 'SuiNetwork'
Set-TextCell 37 3 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextCell 37 4 '0.986'
Set-TextCell 37 5 '  +17.73%  '

# Row 38
Set-TextCell 38 2 'Stacks'
Set-TextCell 38 3 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell 38 4 '1.72'
Set-TextCell 38 5 '  +6.48%  '

# Row 39
Set-TextCell 39 4 '341.51'
Set-TextCell 39 5 '  +10.71%  '

# Row 40
Set-TextCell 40 4 '4.17'
Set-TextCell 40 5 '  +7.01%  '

# Row 41
Set-TextCell 41 4 '38.46'
Set-TextCell 41 5 '  +1.22%  '

# Row 42
Set-TextCell 42 4 '5.33'
Set-TextCell 42 5 '  +6.69%  '

# Row 43
Set-TextCell 43 4 '0.0582'
Set-TextCell 43 5 '  +6.61%  '

# Row 44
Set-TextCell 44 2 'EnergySwap'
Set-TextCell 44 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 44 4 '20.73'
Set-TextCell 44 5 '  +5.66%  '

# Row 45
Set-TextCell 45 2 'Mantle'
Set-TextCell 45 3 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 45 4 '0.632'
Set-TextCell 45 5 '  +5.21%  '

# Row 46
Set-TextCell 46 2 'InjectiveProtocol'
Set-TextCell 46 3 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 46 4 '21.15'
Set-TextCell 46 5 '  +5.90%  '

# Row 47
Set-TextCell 47 4 '135.75'
Set-TextCell 47 5 '  +0.59%  '

# Row 48
Set-TextCell 48 4 '0.0254'
Set-TextCell 48 5 '  +5.55%  '

# Row 49
Set-TextCell 49 4 '0.1000'
Set-TextCell 49 5 '  +1.29%  '

# Row 50
Set-TextCell 50 4 '0.997'
Set-TextCell 50 5 '  -0.13%  '

# Row 51
Set-TextCell 51 4 '2.099.98'
Set-TextCell 51 5 '  +4.26%  '
